$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextCell 'D2' '42.011.13'
Set-TextCell 'E2' '  -0.98%  '
Set-TextCell 'D3' '2.243.23'
Set-TextCell 'E3' '  -1.62%  '
Set-TextCell 'E4' '  -0.03%  '
Set-TextCell 'D5' '305.87'
Set-TextCell 'E5' '  -0.45%  '
Set-TextCell 'D6' '96.25'
Set-TextCell 'E6' '  -1.29%  '
Set-TextCell 'E7' '  -1.59%  '
Set-TextCell 'E8' '  +0.00%  '
Set-TextCell 'D9' '0.485'
Set-TextCell 'E9' '  -1.60%  '
Set-TextCell 'D10' '34.62'
Set-TextCell 'E10' '  -2.88%  '
Set-TextCell 'D11' '0.0804'
Set-TextCell 'E11' '  +0.84%  '
Set-TextCell 'E12' '  +0.36%  '
Set-TextCell 'D13' '6.76'
Set-TextCell 'E13' '  +0.81%  '
Set-TextCell 'D14' '2.594.85'
Set-TextCell 'E14' '  -1.50%  '
Set-TextCell 'D15' '14.39'
Set-TextCell 'E15' '  -0.65%  '
Set-TextCell 'D16' '2.246.39'
Set-TextCell 'E16' '  -1.36%  '
Set-TextCell 'D17' '0.775'
Set-TextCell 'E17' '  -3.12%  '
Set-TextCell 'D18' '41.862.27'
Set-TextCell 'E18' '  -1.11%  '
Set-TextCell 'D19' '12.09'
Set-TextCell 'E19' '  -3.79%  '
Set-TextCell 'D20' '0.0₃0897'
Set-TextCell 'E20' '  -1.55%  '
Set-TextCell 'D21' '5.89'
Set-TextCell 'E21' '  -1.41%  '
Set-TextCell 'D22' '67.00'
Set-TextCell 'E22' '  -1.01%  '
Set-TextCell 'D23' '234.92'
Set-TextCell 'E23' '  -2.42%  '
Set-TextCell 'D24' '2.56'
Set-TextCell 'E24' '  -1.75%  '
Set-TextCell 'E25' '  -0.67%  '
Set-TextCell 'E26' '  -0.04%  '
Set-TextCell 'D27' '23.20'
Set-TextCell 'E27' '  -2.90%  '
Set-TextCell 'D28' '37.47'
Set-TextCell 'E28' '  -0.98%  '
Set-TextCell 'B29' 'Toncoin'
Set-TextCell 'C29' 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextCell 'D29' '2.11'
Set-TextCell 'E29' '  +0.60%  '
Set-TextCell 'B30' 'Cosmos'
Set-TextCell 'C30' 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextCell 'D30' '9.44'
Set-TextCell 'E30' '  -0.75%  '
Set-TextCell 'D31' '165.12'
Set-TextCell 'E31' '  +3.30%  '
Set-TextCell 'E32' '  -0.01%  '
Set-TextCell 'D33' '5.13'
Set-TextCell 'E33' '  -2.66%  '
Set-TextCell 'D34' '3.07'
Set-TextCell 'E34' '  -2.50%  '
Set-TextCell 'D35' '17.44'
Set-TextCell 'E35' '  +2.31%  '
Set-TextCell 'D36' '0.0716'
Set-TextCell 'E36' '  -3.36%  '
Set-TextCell 'E37' '  -0.47%  '
Set-TextCell 'E38' '  -0.87%  '
Set-TextCell 'D39' '0.102'
Set-TextCell 'E39' '  -3.56%  '
Set-TextCell 'D40' '1.78'
Set-TextCell 'E40' '  -3.45%  '
Set-TextCell 'D41' '4.03'
Set-TextCell 'E41' '  -2.47%  '
Set-TextCell 'D42' '1.937.71'
Set-TextCell 'E42' '  -3.15%  '
Set-TextCell 'D43' '0.0279'
Set-TextCell 'E43' '  -2.26%  '
Set-TextCell 'B44' 'ApeXProtocol'
Set-TextCell 'C44' 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
Set-TextCell 'D44' '2.18'
Set-TextCell 'E44' '  -9.94%  '
Set-TextCell 'B45' 'EnergySwap'
Set-TextCell 'C45' 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextCell 'D45' '18.44'
Set-TextCell 'E45' '  -2.85%  '
Set-TextCell 'D46' '2.88'
Set-TextCell 'E46' '  -3.56%  '
Set-TextCell 'D47' '9.63'
Set-TextCell 'E47' '  -3.94%  '
Set-TextCell 'D48' '53.33'
Set-TextCell 'E48' '  +0.62%  '
Set-TextCell 'D49' '2.466.61'
Set-TextCell 'E49' '  -1.46%  '
Set-TextCell 'D50' '70.97'
Set-TextCell 'E50' '  -1.62%  '
Set-TextCell 'D51' '90.80'
Set-TextCell 'E51' '  -1.37%  '
